$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Productdata")
$ws.Range("C2").Value = 0
$ws.Range("E2").Value = 192.53025
$ws.Range("C3").Value = 1899
$ws.Range("E3").Value = 69.84258749999999
$ws.Range("C4").Value = 467
$ws.Range("E4").Value = 16.99281944444444
$ws.Range("C5").Value = 0
$ws.Range("E5").Value = 24.33163958333333
$ws.Range("C6").Value = 967
$ws.Range("E6").Value = 34.17756805555555
$ws.Range("C7").Value = 274
$ws.Range("E7").Value = 10.32906666666667
$ws.Range("C8").Value = 87
$ws.Range("E8").Value = 3.395127777777777
$ws.Range("C9").Value = 1337
$ws.Range("E9").Value = 50.39046736111111
$ws.Range("C10").Value = 958
$ws.Range("E10").Value = 36.08999999999999
$ws.Range("C11").Value = 2855
$ws.Range("E11").Value = 107.5345375
$ws.Range("C12").Value = 669
$ws.Range("E12").Value = 25.1899375
$ws.Range("C13").Value = 0
$ws.Range("E13").Value = 184.4091180555555
$ws.Range("C14").Value = 0
$ws.Range("E14").Value = 68.33114583333332
$ws.Range("C15").Value = 0
$ws.Range("E15").Value = 16.38200347222222
$ws.Range("C16").Value = 0
$ws.Range("E16").Value = 23.33592083333333
$ws.Range("C17").Value = 0
$ws.Range("E17").Value = 33.50895625
$ws.Range("C18").Value = 0
$ws.Range("E18").Value = 10.12695
$ws.Range("C19").Value = 0
$ws.Range("E19").Value = 3.329402777777777
$ws.Range("C20").Value = 0
$ws.Range("E20").Value = 43.76940972222222
$ws.Range("C21").Value = 0
$ws.Range("E21").Value = 46.46547222222221
$ws.Range("C22").Value = 0
$ws.Range("E22").Value = 57.51599999999999
$ws.Range("C23").Value = 0
$ws.Range("E23").Value = 176.9748680555555

$ws = $wb.Worksheets.Item("ForecastedAverageDemand")
$ws.Range("C2").Value = 598
$ws.Range("D2").Value = 153
$ws.Range("G2").Value = 94
$ws.Range("H2").Value = 28
$ws.Range("I2").Value = 418
$ws.Range("J2").Value = 305
$ws.Range("K2").Value = 902
$ws.Range("L2").Value = 206
$ws.Range("C3").Value = 596
$ws.Range("D3").Value = 144
$ws.Range("G3").Value = 77
$ws.Range("H3").Value = 30
$ws.Range("I3").Value = 423
$ws.Range("J3").Value = 297
$ws.Range("K3").Value = 903
$ws.Range("L3").Value = 211
$ws.Range("C4").Value = 598
$ws.Range("D4").Value = 144
$ws.Range("F4").Value = 305
$ws.Range("G4").Value = 88
$ws.Range("H4").Value = 25
$ws.Range("J4").Value = 302
$ws.Range("K4").Value = 889
$ws.Range("L4").Value = 214
$ws.Range("C5").Value = 596
$ws.Range("D5").Value = 140
$ws.Range("G5").Value = 92
$ws.Range("H5").Value = 29
$ws.Range("I5").Value = 420
$ws.Range("J5").Value = 304
$ws.Range("K5").Value = 902
$ws.Range("C6").Value = 604
$ws.Range("D6").Value = 143
$ws.Range("F6").Value = 298
$ws.Range("G6").Value = 96
$ws.Range("H6").Value = 36
$ws.Range("I6").Value = 417
$ws.Range("J6").Value = 302
$ws.Range("K6").Value = 903
$ws.Range("L6").Value = 216
$ws.Range("C7").Value = 610
$ws.Range("D7").Value = 150
$ws.Range("F7").Value = 296
$ws.Range("G7").Value = 90
$ws.Range("H7").Value = 29
$ws.Range("I7").Value = 417
$ws.Range("J7").Value = 304
$ws.Range("K7").Value = 899
$ws.Range("L7").Value = 210
$ws.Range("C8").Value = 593
$ws.Range("D8").Value = 144
$ws.Range("F8").Value = 301
$ws.Range("G8").Value = 91
$ws.Range("H8").Value = 34
$ws.Range("I8").Value = 419
$ws.Range("J8").Value = 292
$ws.Range("L8").Value = 207
$ws.Range("C9").Value = 599
$ws.Range("D9").Value = 147
$ws.Range("F9").Value = 296
$ws.Range("G9").Value = 96
$ws.Range("H9").Value = 28
$ws.Range("I9").Value = 418
$ws.Range("K9").Value = 907

$ws = $wb.Worksheets.Item("ForcastedStandardDeviation")
$ws.Range("C2").Value = 14.95
$ws.Range("D2").Value = 3.824999999999999
$ws.Range("G2").Value = 2.35
$ws.Range("H2").Value = 0.6999999999999998
$ws.Range("I2").Value = 10.45
$ws.Range("J2").Value = 7.624999999999998
$ws.Range("K2").Value = 22.54999999999999
$ws.Range("L2").Value = 5.149999999999999
$ws.Range("C3").Value = 28.30999999999999
$ws.Range("D3").Value = 6.839999999999998
$ws.Range("G3").Value = 3.657499999999999
$ws.Range("H3").Value = 1.425
$ws.Range("I3").Value = 20.09249999999999
$ws.Range("J3").Value = 14.1075
$ws.Range("K3").Value = 42.89249999999999
$ws.Range("L3").Value = 10.0225
$ws.Range("C4").Value = 40.51449999999998
$ws.Range("D4").Value = 9.755999999999997
$ws.Range("F4").Value = 20.66374999999999
$ws.Range("G4").Value = 5.961999999999998
$ws.Range("H4").Value = 1.693749999999999
$ws.Range("J4").Value = 20.46049999999999
$ws.Range("K4").Value = 60.22974999999998
$ws.Range("L4").Value = 14.49849999999999
$ws.Range("C5").Value = 51.2411
$ws.Range("D5").Value = 12.0365
$ws.Range("G5").Value = 7.9097
$ws.Range("H5").Value = 2.493275
$ws.Range("I5").Value = 36.1095
$ws.Range("J5").Value = 26.1364
$ws.Range("K5").Value = 77.54944999999999
$ws.Range("C6").Value = 61.83600999999999
$ws.Range("D6").Value = 14.6399825
$ws.Range("F6").Value = 30.508495
$ws.Range("G6").Value = 9.828239999999997
$ws.Range("H6").Value = 3.685589999999999
$ws.Range("I6").Value = 42.69141749999999
$ws.Range("J6").Value = 30.91800499999999
$ws.Range("K6").Value = 92.44688249999999
$ws.Range("L6").Value = 22.11354
$ws.Range("C7").Value = 71.4552475
$ws.Range("D7").Value = 17.5709625
$ws.Range("F7").Value = 34.67336599999999
$ws.Range("G7").Value = 10.5425775
$ws.Range("H7").Value = 3.397052749999999
$ws.Range("I7").Value = 48.84727574999999
$ws.Range("J7").Value = 35.610484
$ws.Range("K7").Value = 105.30863525
$ws.Range("L7").Value = 24.5993475
$ws.Range("C8").Value = 77.34248457499999
$ws.Range("D8").Value = 18.7813116
$ws.Range("F8").Value = 39.25815827499999
$ws.Range("G8").Value = 11.868745525
$ws.Range("H8").Value = 4.434476349999999
$ws.Range("I8").Value = 54.64839972499998
$ws.Range("J8").Value = 38.08432629999999
$ws.Range("L8").Value = 26.99813542499999
$ws.Range("C9").Value = 85.28753530249999
$ws.Range("D9").Value = 20.9303300325
$ws.Range("F9").Value = 42.14542646
$ws.Range("G9").Value = 13.66878696
$ws.Range("H9").Value = 3.986729529999999
$ws.Range("I9").Value = 59.51617655499999
$ws.Range("K9").Value = 129.1415601325

$ws = $wb.Worksheets.Item("Capacity")
$ws.Range("B2").Value = 48618.75
$ws.Range("B3").Value = 23970
$ws.Range("B4").Value = 7281.25
$ws.Range("B5").Value = 2096.25
$ws.Range("B6").Value = 12005
$ws.Range("B7").Value = 4525
$ws.Range("B8").Value = 1195
$ws.Range("B9").Value = 8382.5
$ws.Range("B10").Value = 15037.5
$ws.Range("B11").Value = 45037.5
$ws.Range("B12").Value = 4192.5
$ws.Range("B13").Value = 64825
$ws.Range("B14").Value = 5992.5
$ws.Range("B15").Value = 1456.25
$ws.Range("B16").Value = 6288.75
$ws.Range("B17").Value = 12005
$ws.Range("B18").Value = 1810
$ws.Range("B19").Value = 597.5
$ws.Range("B20").Value = 149781.25
$ws.Range("B21").Value = 149781.25
$ws.Range("B22").Value = 119825
$ws.Range("B23").Value = 149781.25

$ws = $wb.Worksheets.Item("ProcessingTime")
$ws.Range("B2").Value = 3
$ws.Range("C3").Value = 4
$ws.Range("D4").Value = 5
$ws.Range("E5").Value = 1
$ws.Range("F6").Value = 4
$ws.Range("G7").Value = 5
$ws.Range("H8").Value = 4
$ws.Range("I9").Value = 2
$ws.Range("J10").Value = 5
$ws.Range("K11").Value = 5
$ws.Range("L12").Value = 2
$ws.Range("M13").Value = 4
$ws.Range("N14").Value = 1
$ws.Range("O15").Value = 1
$ws.Range("Q17").Value = 4
$ws.Range("R18").Value = 2
$ws.Range("U21").Value = 5
$ws.Range("V22").Value = 4
$ws.Range("W23").Value = 5
